# Add three new rows (124-126) of uploaded/collected file records to each
# of the four log sheets, extending the data range from A1:I123 to A1:I126.

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# New timestamps (Excel date serials) shared by the three new rows on every sheet.
$newTimes = @(45910.4328125, 45911.43440972222, 45912.43680555555)

# Per-sheet, per-row data for columns B..I.
# Each sheet has an array of 3 rows (for r=124,125,126), each row is:
# @(B, C, D, E, F, G, H, I)
$g1 = [double]"7.598631275147109e+23"
$g2 = [double]"5.68432987514711e+23"
$g3 = [double]"5.68631262647114e+23"
$g4 = [double]"9.85046333984776e+23"

$sheetData = @{
    "DE_LFT_#1" = @(
        @("0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x01,0x08", "0x14", 380, $g1, 264, 14),
        @("0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x01,0x08", "0x14", 380, $g1, 264, 14),
        @("0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x01,0x08", "0x14", 380, $g1, 264, 14)
    )
    "DE_LFT_#2" = @(
        @("0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x01,0x10", "0xe", 380, $g2, 272, 14),
        @("0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x01,0x0C", "0xe", 380, $g2, 268, 14),
        @("0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x01,0x0C", "0xe", 380, $g2, 268, 14)
    )
    "DE_PLT_#1" = @(
        @("0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x70", "0x7", 130, $g3, 112, 7),
        @("0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x70", "0x7", 130, $g3, 112, 7),
        @("0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x70", "0x7", 130, $g3, 112, 7)
    )
    "DE_PLT_#2" = @(
        @("0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x6E", "0x3", 130, $g4, 110, 3),
        @("0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x6E", "0x3", 130, $g4, 110, 3),
        @("0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x6D", "0x3", 130, $g4, 109, 3)
    )
}

foreach ($ws in $wb.Worksheets) {
    $rowsForSheet = $sheetData[$ws.Name]
    if (-not $rowsForSheet) { continue }

    for ($i = 0; $i -lt 3; $i++) {
        $rowNum = 124 + $i
        $vals = $rowsForSheet[$i]

        $cellA = $ws.Cells.Item($rowNum, 1)
        $cellA.NumberFormat = $dateFormat
        $cellA.Value = $newTimes[$i]

        $ws.Cells.Item($rowNum, 2).Value = $vals[0]
        $ws.Cells.Item($rowNum, 3).Value = $vals[1]
        $ws.Cells.Item($rowNum, 4).Value = $vals[2]
        $ws.Cells.Item($rowNum, 5).Value = $vals[3]
        $ws.Cells.Item($rowNum, 6).Value = $vals[4]
        $ws.Cells.Item($rowNum, 7).Value = $vals[5]
        $ws.Cells.Item($rowNum, 8).Value = $vals[6]
        $ws.Cells.Item($rowNum, 9).Value = $vals[7]
    }
}
